$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.962.37"
$ws.Range("E2").Value = "  +3.28%  "

$ws.Range("D3").Value = "1.724.90"
$ws.Range("E3").Value = "  +2.95%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.57"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.01"
$ws.Range("E8").Value = "  +13.04%  "

$ws.Range("E9").Value = "  +3.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0633"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("D12").Value = "1.969.07"
$ws.Range("E12").Value = "  +2.96%  "

$ws.Range("D13").Value = "1.721.93"
$ws.Range("E13").Value = "  +2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.27"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  +5.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.91"
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("D17").Value = "27.909.02"
$ws.Range("E17").Value = "  +3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.39"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("D19").Value = "0.0₃0755"
$ws.Range("E19").Value = "  +2.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -3.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("E22").Value = "  +4.09%  "

$ws.Range("E23").Value = "  +4.85%  "

$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.43"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("E26").Value = "  +4.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.81"
$ws.Range("E27").Value = "  +3.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +1.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  +3.02%  "

$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("E33").Value = "  +3.68%  "

$ws.Range("D34").Value = "1.491.77"
$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.613"
$ws.Range("E36").Value = "  +3.49%  "

$ws.Range("E37").Value = "  +4.99%  "

$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.44"
$ws.Range("E41").Value = "  +5.67%  "

$ws.Range("E42").Value = "  +5.88%  "

$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").Value = "1.873.18"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.792"
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("E47").Value = "  +12.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.09"
$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.27"
$ws.Range("E50").Value = "  +3.85%  "

$ws.Range("E51").Value = "  +1.49%  "
